# edit.ps1 -- "ajout dans outil et source"
#
# Applies, via the Word object model, the same edits the author made by
# hand in Word:
#   1. expands the opening sentence of "Présentation du projet"
#   2. corrects "100 ans" to "105 ans" (and drags the _GoBack marker
#      along, since that's where the author's cursor ended up last)
#   3. retitles "Outil de travail" -> "Outil et source"
#   4. splits "...travaillé en équipe. Nous avons utilisé les
#      logiciels :" into its own paragraph
#   5. expands the "Sublime Text" bullet and appends a new
#      "Microsoft Word" bullet
#
# Word constants used below: wdFindContinue = 1, wdReplaceAll = 2

$d = $word.ActiveDocument

function New-OpenXmlParagraphFragment([string]$innerXml) {
    return '<?xml version="1.0"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>' + $innerXml + '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
}

# ---------------------------------------------------------------------
# 1) "Nous avons choisi de modéliser un automate cellulaire" ->
#    "Nous avons choisi de modéliser, dans le langage de programmation
#     Python, un automate cellulaire"
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "Nous avons choisi de modéliser un automate cellulaire",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "Nous avons choisi de modéliser, dans le langage de programmation Python, un automate cellulaire",
    2)

# ---------------------------------------------------------------------
# 2) "(0 à 100 ans)" -> "(0 à 105 ans)", and move the "_GoBack" bookmark
#    so it sits right after the newly typed "5" (its last-edit spot).
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    " (0 à 100 ans)",
    $false, $false, $false, $false, $false, $true, 1, $false,
    " (0 à 105 ans)",
    2)

$d.Bookmarks("_GoBack").Delete()
$fiveSpot = $d.Content
$fiveSpot.Find.Execute(" (0 à 105")
$markPoint = $d.Range($fiveSpot.End, $fiveSpot.End)
$d.Bookmarks.Add("_GoBack", $markPoint)

# ---------------------------------------------------------------------
# 3) Title "Outil de travail" -> "Outil et source"
#    Built with InsertXML so the existing <w:lastRenderedPageBreak/>
#    stays attached to the first run, exactly like the diff shows.
# ---------------------------------------------------------------------
$titleRange = $d.Content
$titleRange.Find.Execute("Outil de travail")
$titleStart = $d.Range($titleRange.Start, $titleRange.Start)
$titleFragment = '<w:r><w:lastRenderedPageBreak/><w:t xml:space="preserve">Outil </w:t></w:r><w:r><w:t>et source</w:t></w:r>'
$titleStart.InsertXML((New-OpenXmlParagraphFragment $titleFragment))
$d.Content.Find.Execute(
    "Outil de travail",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "",
    2)

# ---------------------------------------------------------------------
# 4) Split "... travaillé en équipe. Nous avons utilisé les
#    logiciels :" into two paragraphs (this also swallows the old
#    "_GoBack" bookmark that used to live right in the middle of it).
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "travaillé en équipe. Nous avons utilisé les logiciels" + [char]160 + ":",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "travaillé en équipe.^pNous avons utilisé les logiciels :",
    2)

# ---------------------------------------------------------------------
# 5) Expand the "Sublime Text" bullet (keeping/augmenting the spell-
#    check proofErr markers around "Text" and "bash") and append a new
#    "Microsoft Word" bullet right after it.
# ---------------------------------------------------------------------
$bulletRange = $d.Content
$bulletRange.Find.Execute("Sublime Text qui est open source, ")
$bulletStart = $d.Range($bulletRange.Start, $bulletRange.Start)
$bulletFragment = '<w:r><w:t xml:space="preserve">Sublime </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/><w:r><w:t>Text</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve"> : </w:t></w:r>' +
    '<w:r><w:t>est open source</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> qui est exécutable sous le </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/><w:r><w:t>bash</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve"> Linux et prend en compte le langage Python </w:t></w:r>' +
    '<w:r><w:t>(tout comme</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> 44</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> langage</w:t></w:r>' +
    '<w:r><w:t>s</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> au totale)</w:t></w:r>'
$bulletStart.InsertXML((New-OpenXmlParagraphFragment $bulletFragment))
$d.Content.Find.Execute(
    "Sublime Text qui est open source, ",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "",
    2)

$d.Content.Find.Execute(
    " au totale)",
    $false, $false, $false, $false, $false, $true, 1, $false,
    " au totale)^pMicrosoft Word : pour toute la rédaction du présent dossier",
    2)
